# Apply "updated SIL+AUC database outputs" edit (area 20-22, 2024 escapement prelim-inseason)
# to the EscEstSppHeader sheet: revised survey stats for several existing rows (27-29, 31,
# 36, 40-44) plus four brand-new rows (45-48) pushing the table from 44 to 48 data rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D27").Value = (Get-Date -Year 2024 -Month 7 -Day 28 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E27").Value = (Get-Date -Year 2024 -Month 11 -Day 30 -Hour 0 -Minute 0 -Second 0)
$ws.Range("F27").Value = 8
$ws.Range("G27").Value = 10
$ws.Range("H27").Value = 8
$ws.Range("I27").Value = 10
$ws.Range("J27").Value = 125
$ws.Range("K27").Value = 2065028.5
$ws.Range("L27").Value = 206503
$ws.Range("M27").Value = 258129
$ws.Range("Q27").Value = 2393962.84358531
$ws.Range("R27").Value = 239396
$ws.Range("S27").Value = 299245
$ws.Range("W27").Value = 102250
$ws.Range("X27").Value = 113278
$ws.Range("AD27").ClearContents()
$ws.Range("AM27").Value = 'Sept 20 % pop reduced as chum were starting to arrive and may have missed some in the upper section which was not swum this day.'

$ws.Range("D28").Value = (Get-Date -Year 2024 -Month 8 -Day 15 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J28").Value = 92
$ws.Range("K28").Value = 233865
$ws.Range("L28").Value = 9355
$ws.Range("M28").Value = 15591
$ws.Range("N28").Value = 4421.5
$ws.Range("O28").Value = 177
$ws.Range("P28").Value = 295
$ws.Range("Q28").Value = 303095.482360849
$ws.Range("R28").Value = 12124
$ws.Range("S28").Value = 20206
$ws.Range("T28").Value = 5449.96066888726
$ws.Range("U28").Value = 218
$ws.Range("V28").Value = 363
$ws.Range("X28").Value = 9007
$ws.Range("Z28").Value = 364
$ws.Range("AM28").Value = 'Sept 4 and 20 % pop reduced as chum were starting to arrive and may have missed some in the upper section which was not swum on either day. Oct 24 % pop reduced further as only swam the lower; see accompanying prelim esc Rmarkdown file for details.'

$ws.Range("D29").Value = (Get-Date -Year 2024 -Month 8 -Day 15 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E29").Value = (Get-Date -Year 2024 -Month 12 -Day 31 -Hour 0 -Minute 0 -Second 0)
$ws.Range("F29").Value = 25
$ws.Range("G29").Value = 35
$ws.Range("H29").Value = 25
$ws.Range("I29").Value = 35
$ws.Range("J29").Value = 138
$ws.Range("K29").Value = 66034.5
$ws.Range("L29").Value = 1887
$ws.Range("M29").Value = 2641
$ws.Range("N29").Value = 7393.5
$ws.Range("O29").Value = 211
$ws.Range("P29").Value = 296
$ws.Range("Q29").Value = 85972.6128382703
$ws.Range("R29").Value = 2456
$ws.Range("S29").Value = 3439
$ws.Range("T29").Value = 9682.72968490306
$ws.Range("U29").Value = 277
$ws.Range("V29").Value = 387
$ws.Range("X29").Value = 1793
$ws.Range("Z29").Value = 203
$ws.Range("AM29").Value = 'Oct 24 % pop reduced further as only swam the lower; see accompanying prelim esc Rmarkdown file for details. Nov 7 % pop reduced as did not swim tributaries where coho may be heading to this time of year.'

$ws.Range("D31").Value = (Get-Date -Year 2024 -Month 7 -Day 1 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E31").Value = (Get-Date -Year 2024 -Month 11 -Day 10 -Hour 0 -Minute 0 -Second 0)
$ws.Range("F31").Value = 30
$ws.Range("G31").Value = 40
$ws.Range("H31").Value = 30
$ws.Range("I31").Value = 40
$ws.Range("J31").Value = 132
$ws.Range("K31").Value = 11274
$ws.Range("L31").Value = 282
$ws.Range("M31").Value = 376
$ws.Range("P31").Value = 0
$ws.Range("Q31").Value = 14560.4702426181
$ws.Range("R31").Value = 364
$ws.Range("S31").Value = 485
$ws.Range("X31").Value = 281
$ws.Range("AA31").Value = 1
$ws.Range("AB31").Value = 2
$ws.Range("AC31").Value = 2
$ws.Range("AD31").Value = 353
$ws.Range("AM31").Value = '% pop lowered for almost all surveys as did not survey Parker Creek/Tuck Lake, and surveys were not timed for SK. Oct 24 % pop reduced further as only swam the lower; see accompanying prelim esc RMarkdown file for details.'

$ws.Range("Q36").Value = 78.1249968955915
$ws.Range("R36").Value = 3
$ws.Range("S36").Value = 5
$ws.Range("X36").Value = 2
$ws.Range("AD36").Value = 2
$ws.Range("AM36").Value = 'One swim Nov 1. Swam marker 24 to lake, as well as side channel and Borrow Pits. The percent population was adjusted to 60% for three reasons: this survey occurred late for Chinook , the creek recently flooded prior to the swim, and to account for the swampy and highly vegetated channels that were not surveyed up above. OE of 80% used since swim conditions were good; water level was normal and visibility was good to slightly turbid in some sections. Type-6 due to only one survey and well after peak spawning. Used PL+D expanded.'

$ws.Range("A40").Value = '930-071700-37600-00000-0000-0000-000-000-000-000-000-000'
$ws.Range("K40").Value = 14250
$ws.Range("L40").Value = 570
$ws.Range("M40").Value = 950
$ws.Range("Q40").Value = 14250
$ws.Range("R40").Value = 570
$ws.Range("S40").Value = 950
$ws.Range("W40").Value = 380
$ws.Range("X40").Value = 380
$ws.Range("AA40").Value = 2
$ws.Range("AD40").ClearContents()

$ws.Range("A41").Value = '930-071700-37600-00000-0000-0000-000-000-000-000-000-000'

$ws.Range("A42").Value = '930-071700-37600-00000-0000-0000-000-000-000-000-000-000'
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = 0
$ws.Range("N42").Value = 0
$ws.Range("O42").Value = 0
$ws.Range("P42").Value = 0
$ws.Range("Q42").Value = 0
$ws.Range("R42").Value = 0
$ws.Range("S42").Value = 0
$ws.Range("T42").Value = 0
$ws.Range("U42").Value = 0
$ws.Range("V42").Value = 0
$ws.Range("W42").Value = 0
$ws.Range("X42").Value = 0
$ws.Range("Y42").Value = 0
$ws.Range("Z42").Value = 0

$ws.Range("A43").Value = '930-071700-37600-00000-0000-0000-000-000-000-000-000-000'

$ws.Range("A44").Value = '930-071700-44200-00000-0000-0000-000-000-000-000-000-000'
$ws.Range("C44").Value = 'CM'
$ws.Range("J44").Value = 75
$ws.Range("K44").Value = 13412
$ws.Range("L44").Value = 536
$ws.Range("M44").Value = 894
$ws.Range("Q44").Value = 14902.2226169963
$ws.Range("R44").Value = 596
$ws.Range("S44").Value = 993
$ws.Range("W44").Value = 213
$ws.Range("X44").Value = 216

$ws.Range("A45").Value = '930-071700-44200-00000-0000-0000-000-000-000-000-000-000'
$ws.Range("B45").Value = 2024
$ws.Range("C45").Value = 'CN'
$ws.Range("D45").Value = (Get-Date -Year 2024 -Month 9 -Day 1 -Hour 0 -Minute 0 -Second 0)
$ws.Range("D45").NumberFormat = "mm-dd-yy"
$ws.Range("E45").Value = (Get-Date -Year 2024 -Month 11 -Day 15 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E45").NumberFormat = "mm-dd-yy"
$ws.Range("F45").Value = 15
$ws.Range("G45").Value = 25
$ws.Range("H45").Value = 15
$ws.Range("I45").Value = 25
$ws.Range("J45").Value = 75
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = 0
$ws.Range("N45").Value = 0
$ws.Range("O45").Value = 0
$ws.Range("P45").Value = 0
$ws.Range("Q45").Value = 0
$ws.Range("R45").Value = 0
$ws.Range("S45").Value = 0
$ws.Range("T45").Value = 0
$ws.Range("U45").Value = 0
$ws.Range("V45").Value = 0
$ws.Range("W45").Value = 0
$ws.Range("X45").Value = 0
$ws.Range("Y45").Value = 0
$ws.Range("Z45").Value = 0
$ws.Range("AA45").Value = 0
$ws.Range("AB45").Value = 0
$ws.Range("AC45").Value = 0
$ws.Range("AD45").Value = 0
$ws.Range("AE45").Value = 0
$ws.Range("AF45").Value = 0
$ws.Range("AG45").Value = 0
$ws.Range("AH45").Value = $false
$ws.Range("AI45").Value = $false
$ws.Range("AJ45").Value = 0
$ws.Range("AK45").Value = 0
$ws.Range("AL45").Value = -1
$ws.Range("AN45").Value = $false

$ws.Range("A46").Value = '930-071700-44200-00000-0000-0000-000-000-000-000-000-000'
$ws.Range("B46").Value = 2024
$ws.Range("C46").Value = 'CO'
$ws.Range("D46").Value = (Get-Date -Year 2024 -Month 9 -Day 1 -Hour 0 -Minute 0 -Second 0)
$ws.Range("D46").NumberFormat = "mm-dd-yy"
$ws.Range("E46").Value = (Get-Date -Year 2024 -Month 11 -Day 15 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E46").NumberFormat = "mm-dd-yy"
$ws.Range("F46").Value = 15
$ws.Range("G46").Value = 25
$ws.Range("H46").Value = 15
$ws.Range("I46").Value = 25
$ws.Range("J46").Value = 75
$ws.Range("K46").Value = 396
$ws.Range("L46").Value = 16
$ws.Range("M46").Value = 26
$ws.Range("N46").Value = 88
$ws.Range("O46").Value = 4
$ws.Range("P46").Value = 6
$ws.Range("Q46").Value = 440.00001165602
$ws.Range("R46").Value = 18
$ws.Range("S46").Value = 29
$ws.Range("T46").Value = 97.7777803680044
$ws.Range("U46").Value = 4
$ws.Range("V46").Value = 7
$ws.Range("W46").Value = 36
$ws.Range("X46").Value = 40
$ws.Range("Y46").Value = 8
$ws.Range("Z46").Value = 9
$ws.Range("AA46").Value = 0
$ws.Range("AB46").Value = 0
$ws.Range("AC46").Value = 0
$ws.Range("AD46").Value = 0
$ws.Range("AE46").Value = 0
$ws.Range("AF46").Value = 0
$ws.Range("AG46").Value = 0
$ws.Range("AH46").Value = $false
$ws.Range("AI46").Value = $false
$ws.Range("AJ46").Value = 0
$ws.Range("AK46").Value = 0
$ws.Range("AL46").Value = -1
$ws.Range("AN46").Value = $false

$ws.Range("A47").Value = '930-071700-44200-00000-0000-0000-000-000-000-000-000-000'
$ws.Range("B47").Value = 2024
$ws.Range("C47").Value = 'SK'
$ws.Range("D47").Value = (Get-Date -Year 2024 -Month 9 -Day 1 -Hour 0 -Minute 0 -Second 0)
$ws.Range("D47").NumberFormat = "mm-dd-yy"
$ws.Range("E47").Value = (Get-Date -Year 2024 -Month 11 -Day 15 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E47").NumberFormat = "mm-dd-yy"
$ws.Range("F47").Value = 15
$ws.Range("G47").Value = 25
$ws.Range("H47").Value = 15
$ws.Range("I47").Value = 25
$ws.Range("J47").Value = 75
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = 0
$ws.Range("N47").Value = 0
$ws.Range("O47").Value = 0
$ws.Range("P47").Value = 0
$ws.Range("Q47").Value = 0
$ws.Range("R47").Value = 0
$ws.Range("S47").Value = 0
$ws.Range("T47").Value = 0
$ws.Range("U47").Value = 0
$ws.Range("V47").Value = 0
$ws.Range("W47").Value = 0
$ws.Range("X47").Value = 0
$ws.Range("Y47").Value = 0
$ws.Range("Z47").Value = 0
$ws.Range("AA47").Value = 0
$ws.Range("AB47").Value = 0
$ws.Range("AC47").Value = 0
$ws.Range("AD47").Value = 0
$ws.Range("AE47").Value = 0
$ws.Range("AF47").Value = 0
$ws.Range("AG47").Value = 0
$ws.Range("AH47").Value = $false
$ws.Range("AI47").Value = $false
$ws.Range("AJ47").Value = 0
$ws.Range("AK47").Value = 0
$ws.Range("AL47").Value = -1
$ws.Range("AN47").Value = $false

$ws.Range("A48").Value = '930-979400-00000-00000-0000-0000-000-000-000-000-000-000'
$ws.Range("B48").Value = 2024
$ws.Range("C48").Value = 'CN'
$ws.Range("D48").Value = (Get-Date -Year 2024 -Month 9 -Day 1 -Hour 0 -Minute 0 -Second 0)
$ws.Range("D48").NumberFormat = "mm-dd-yy"
$ws.Range("E48").Value = (Get-Date -Year 2024 -Month 11 -Day 15 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E48").NumberFormat = "mm-dd-yy"
$ws.Range("F48").Value = 15
$ws.Range("G48").Value = 25
$ws.Range("H48").Value = 15
$ws.Range("I48").Value = 25
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("M48").Value = 0
$ws.Range("N48").Value = 0
$ws.Range("O48").Value = 0
$ws.Range("P48").Value = 0
$ws.Range("Q48").Value = 0
$ws.Range("R48").Value = 0
$ws.Range("S48").Value = 0
$ws.Range("T48").Value = 0
$ws.Range("U48").Value = 0
$ws.Range("V48").Value = 0
$ws.Range("W48").Value = 0
$ws.Range("X48").Value = 0
$ws.Range("Y48").Value = 0
$ws.Range("Z48").Value = 0
$ws.Range("AA48").Value = 0
$ws.Range("AB48").Value = 0
$ws.Range("AC48").Value = 0
$ws.Range("AD48").Value = 0
$ws.Range("AE48").Value = 0
$ws.Range("AF48").Value = 0
$ws.Range("AG48").Value = 0
$ws.Range("AH48").Value = $false
$ws.Range("AI48").Value = $false
$ws.Range("AJ48").Value = 0
$ws.Range("AK48").Value = 0
$ws.Range("AL48").Value = -1
$ws.Range("AN48").Value = $false

# Rows 32 and 34 had empty placeholder inline-string cells in AM that the
# diff drops entirely; explicitly clear them so they round-trip as truly
# empty cells instead of empty-string cells.
$ws.Range("AM32").ClearContents()
$ws.Range("AM34").ClearContents()

# Update the defined name range to cover the newly added rows (44 -> 48).
$definedName = $wb.Names.Item("EscEstSppHeader")
$definedName.RefersTo = "='EscEstSppHeader'!`$A`$1:`$AN`$48"

